$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first data row (row 2, date 39400 / year 2007) is removed entirely.
# This shifts every subsequent row up by one (rows 3..19 become rows 2..18).
$ws.Rows.Item(2).Delete()

# The y_1_forecast column (E) values are recomputed for the new data window.
$eValues = @(
    0.8212989654785341,
    1.183007486132071,
    1.015842920196763,
    0.9092565586104273,
    1.236730309040235,
    1.029015928490629,
    1.358148715145191,
    1.528208222695326,
    1.634928000057778,
    1.67176973076042,
    1.603287858019664,
    0.8408455317168162,
    -1.875058665585216,
    5.03478667886097,
    2.399708479013141,
    0.8520283695166997,
    0.299857156820571
)

for ($i = 0; $i -lt $eValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 5).Value = $eValues[$i]
}
